# Updates crypto price/volume data per commit "Updated cryptos list on Fri Aug 11 21:17:44 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '29.402.08'
$ws.Range("E2").Value = '  -0.05%  '
$ws.Range("D3").Value = '1.843.99'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'239.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = "'0.9998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'0.07536"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").Value = "'0.2930"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.39%  '
$ws.Range("D10").Value = "'24.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").Value = "'0.07711"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.01%  '
$ws.Range("D12").Value = '1.861.71'
$ws.Range("E12").Value = '  -6.21%  '
$ws.Range("D13").Value = "'4.999"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.30%  '
$ws.Range("D14").Value = "'0.6789"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.77%  '
$ws.Range("D15").Value = "'0.00001039"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.31%  '
$ws.Range("D16").Value = "'83.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '2.136.84'
$ws.Range("E17").Value = '  -5.63%  '
$ws.Range("D18").Value = "'6.162"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").Value = '29.435.84'
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("E21").Value = '  -0.73%  '
$ws.Range("D22").Value = "'0.9998"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").Value = "'7.459"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.60%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = "'157.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  +0.40%  '
$ws.Range("D27").Value = "'8.370"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.66%  '
$ws.Range("D28").Value = "'17.58"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("D29").Value = "'1.458"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("D30").Value = "'1.278"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("D31").Value = "'0.05628"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.73%  '
$ws.Range("D32").Value = "'4.105"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").Value = "'4.030"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.25%  '
$ws.Range("E34").Value = '  -1.48%  '
$ws.Range("D35").Value = "'1.157"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.11%  '
$ws.Range("D36").Value = "'0.7116"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.68%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").Value = '1.246.32'
$ws.Range("D39").Value = "'0.01807"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.13%  '
$ws.Range("D40").Value = "'2.763"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.06%  '
$ws.Range("D41").Value = "'6.319"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.72%  '
$ws.Range("D42").Value = "'0.9004"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.05%  '
$ws.Range("D43").Value = "'0.9996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '
$ws.Range("D44").Value = "'101.86"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.18%  '
$ws.Range("D45").Value = "'65.78"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.77%  '
$ws.Range("D46").Value = "'0.00000000119"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.95%  '
$ws.Range("D47").Value = "'7.088"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.41%  '
$ws.Range("D48").Value = "'0.4002"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = "'1.672"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = "'8.888"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.79%  '
$ws.Range("D51").Value = "'0.1120"
$ws.Range("D51").Style = "Normal"
